# Tidsregistrering - fill in the remaining days of the week (16/3, 20/3 - 24/3)
# that the author had not gotten around to entering yet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# --- Formats to copy from existing rows, so no new style entries are created ---
$dateStyleSrc = $ws.Range("A31")   # numFmtId 14 (date)
$timeStyleSrc = $ws.Range("G5")    # numFmtId 20 (time)

# Date cells (col A) - one per day header row
$dateCells = @("A35","A39","A43","A46","A50","A54")
$dateStyleSrc.Copy()
foreach ($addr in $dateCells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Time cells (cols G/H, and I40 which also carries the time style)
$timeCells = @("G40","H40","I40","G44","H44","G47","H47","G48","H48","G51","H51","G52","H52","G55","H55")
$timeStyleSrc.Copy()
foreach ($addr in $timeCells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# --- Date values (Excel 1900-date-system serials) ---
$ws.Range("A35").Value = 42810   # 16-03-2017
$ws.Range("A39").Value = 42814   # 20-03-2017
$ws.Range("A43").Value = 42815   # 21-03-2017
$ws.Range("A46").Value = 42816   # 22-03-2017
$ws.Range("A50").Value = 42817   # 23-03-2017
$ws.Range("A54").Value = 42818   # 24-03-2017

# --- Activity rows ---
# New shared strings must be created in this order so they land at the same
# sharedStrings.xml indices as the authored workbook: bugfixing, gui - bugfixing, design.
$ws.Range("F48").Value = "bugfixing"
$ws.Range("G48").Value = 0.5
$ws.Range("H48").Value = 0.66666666666666663
$ws.Range("I48").Value = 4

$ws.Range("F51").Value = "gui - bugfixing"
$ws.Range("G51").Value = 0.41666666666666669
$ws.Range("H51").Value = 0.66666666666666663
$ws.Range("I51").Value = 6

$ws.Range("F40").Value = "design"
$ws.Range("G40").Value = 0.33680555555555558
$ws.Range("H40").Value = 0.64583333333333337
$ws.Range("I40").Value = 7

# Remaining rows re-use already-existing shared strings ("test", "gui - bugfixing")
$ws.Range("F44").Value = "test"
$ws.Range("G44").Value = 0.33680555555555558
$ws.Range("H44").Value = 0.64583333333333337
$ws.Range("I44").Value = 7

$ws.Range("F47").Value = "test"
$ws.Range("G47").Value = 0.33680555555555558
$ws.Range("H47").Value = 0.47916666666666669
$ws.Range("I47").Value = 3.5

$ws.Range("F52").Value = "test"
$ws.Range("G52").Value = 0.33680555555555558
$ws.Range("H52").Value = 0.41666666666666669
$ws.Range("I52").Value = 2

$ws.Range("F55").Value = "gui - bugfixing"
$ws.Range("G55").Value = 0.33680555555555558
$ws.Range("H55").Value = 0.5
$ws.Range("I55").Value = 4

# --- Sheet view bookkeeping: scrolled down to the newly filled-in rows ---
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("I32").Select()
